# [DSC-721] Bundle creation from bulk import does not set file format
#
# The "bitstream-metadata" sheet's sample row uses a dc.title value of
# "Test title" for the bitstream metadata. Since bulk import derives the
# bitstream/file format from the file name's extension, the sample title
# needs a realistic file extension (".txt") so the bundle/bitstream format
# detection test data is meaningful.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("bitstream-metadata")

# D2 holds the dc.title value for the sample bitstream row ("Test title").
$ws.Cells.Item(2, 4).Value = "Test title.txt"

# Leave the active selection on that row/sheet, matching the authored edit.
$ws.Activate()
$ws.Range("D3").Select()
